$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "; online test" markers next to the online-test related rows.
$ws.Cells.Item(8, 3).Value = "; online test"
$ws.Cells.Item(9, 3).Value = "; online test"

# Insert a new row (10) for "language" / "English", pushing "test paper" and
# "marksheet" rows down by one.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).RowHeight = 19

$ws.Cells.Item(10, 1).Value = "language"
$ws.Cells.Item(10, 2).Value = "English"
$ws.Cells.Item(10, 3).Value = "; paper test"

# "test paper" row (now row 11): drop the ".pdf" suffix from the value and
# annotate it.
$ws.Cells.Item(11, 2).Value = "testpaper"
$ws.Cells.Item(11, 3).Value = "; paper test"
$ws.Cells.Item(11, 4).Value = "the file type (pdf) will be added by the program"

# "marksheet" row (now row 12): the value mirrors the label, drop the
# ".xlsx" suffix, and annotate it.
$ws.Cells.Item(12, 2).Value = "marksheet"
$ws.Cells.Item(12, 3).Value = "; paper test"
$ws.Cells.Item(12, 4).Value = "the file type (xlsx) will be added by the program"

$ws.Range("A1").Select() | Out-Null
